$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121, shifting existing rows 121.. down by one.
$ws.Rows(121).Insert()

# Populate the newly inserted row 121 with the new record's values.
$ws.Cells.Item(121, 1).Value = 4
$ws.Cells.Item(121, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value = "Los Lagos"
$ws.Cells.Item(121, 4).Value = 44553
$ws.Cells.Item(121, 5).Value = 10
$ws.Cells.Item(121, 6).Value = 100112024
$ws.Cells.Item(121, 7).Value = "Choclo"
$ws.Cells.Item(121, 8).Value = "Dulce o Americano"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 200
$ws.Cells.Item(121, 11).Value = 20000
$ws.Cells.Item(121, 12).Value = 20000
$ws.Cells.Item(121, 13).Value = 20000
$ws.Cells.Item(121, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(121, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(121, 16).Value = 286
$ws.Cells.Item(121, 17).Value = 70
$ws.Cells.Item(121, 18).Value = "Hortaliza"

# Match the date style used by the other rows in column D.
$ws.Cells.Item(121, 4).NumberFormat = $ws.Cells.Item(122, 4).NumberFormat
